# Update the StructureDefinition-level-of-care workbook:
#   - bump Version 5.0.0 -> 6.0.0
#   - bump Date to the new publish timestamp
#   - fill in Publisher ("Alvearie Team")
#   - replace the two "Contact / No display for ContactDetail" rows with a
#     single "Jurisdiction / United States of America" row
#   - update the Elements sheet's top-level Extension row (Short/Definition)
#     to describe "Level Of Care" / "Code for the level of care"

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# --- Version / Date ---------------------------------------------------
$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# --- Collapse the duplicated "Contact" rows (10 & 11) into one row,
#     shifting everything below up by one. ------------------------------
$meta.Rows.Item(11).Delete()

# --- Publisher / Jurisdiction ------------------------------------------
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet: top-level Extension row's Short/Definition --------
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Level Of Care"
$elements.Range("L2").Value = "Code for the level of care"
